# Daily attendance processing - swap the order of "Recorded By" contributors
# in column G: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = [string]$cell.Text
    if ($text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
